# lagt in lyornas area.
# Row 5 ("GIS-data lyornas avstand till rodravslyor") moves from "pabörjat"
# (in progress) to "klar" (done), and its explanatory comment in column C
# (about still missing the area measurement) is removed since the area has
# now been entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-use the formatting already used for other "klar" cells (e.g. B2) so the
# green font style carries over correctly, then set the new value.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B5").Value2 = "klar"

# The "missing area" comment no longer applies.
$ws.Range("C5").ClearContents() | Out-Null

# Reflect the new active selection/cursor position recorded in the workbook.
$ws.Range("C5").Select() | Out-Null
